$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 208.75
$ws.Range("I92").Value = 208.75
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 208.75
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 1039.25
$ws.Range("N92").ClearContents()

$ws.Range("H96").Value = 1577.8182
$ws.Range("I96").Value = 294.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 883.5
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = 489.5
$ws.Range("N96").Value = -17746

$ws.Range("H106").Value = 1743.5
$ws.Range("I106").Value = 1706.8572
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 1706.8572
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -1075.8572
$ws.Range("N106").Value = -3262

$ws.Range("H132").Value = 1042.5667
$ws.Range("I132").Value = 1014.03705
$ws.Range("J132").Value = 1299.3334
$ws.Range("K132").Value = 3042.11115
$ws.Range("L132").Value = 3898.0002
$ws.Range("M132").Value = -512.1111500000002
$ws.Range("N132").Value = -8958.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5818208.5
$ws.Range("I2").Value = 5818208.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5818208.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -5818095.5

$ws.Range("H32").Value = 2138.5615
$ws.Range("I32").Value = 1661.678
$ws.Range("J32").Value = 4148.2856
$ws.Range("K32").Value = 1661.678
$ws.Range("L32").Value = 4148.2856
$ws.Range("M32").Value = -1374.678
$ws.Range("N32").Value = -4722.2856

$ws.Range("H63").Value = 4655
$ws.Range("I63").Value = 4655
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4655
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3969

$ws.Range("H66").Value = 4655
$ws.Range("I66").Value = 4655
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 23275
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -19843

$ws.Range("H74").Value = 1882.2
$ws.Range("I74").Value = 703.1429000000001
$ws.Range("J74").Value = 4633.3335
$ws.Range("K74").Value = 703.1429000000001
$ws.Range("L74").Value = 4633.3335
$ws.Range("M74").Value = 170.8570999999999
$ws.Range("N74").Value = -6381.3335

$ws.Range("H77").Value = 1882.2
$ws.Range("I77").Value = 703.1429000000001
$ws.Range("J77").Value = 4633.3335
$ws.Range("K77").Value = 3515.7145
$ws.Range("L77").Value = 23166.6675
$ws.Range("M77").Value = 852.2855
$ws.Range("N77").Value = -31902.6675

$ws.Range("H116").Value = 5818208.5
$ws.Range("I116").Value = 5818208.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5818208.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -5815914.5

$ws.Range("H132").Value = 2046.9623
$ws.Range("I132").Value = 1806.1777
$ws.Range("J132").Value = 3401.375
$ws.Range("K132").Value = 5418.5331
$ws.Range("L132").Value = 10204.125
$ws.Range("M132").Value = -2888.5331
$ws.Range("N132").Value = -15264.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5818208.5
$ws.Range("I3").Value = 5818208.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5818208.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5818094.5

$ws.Range("H134").Value = 15151.25
$ws.Range("I134").Value = 18885
$ws.Range("J134").Value = 3950
$ws.Range("K134").Value = 56655
$ws.Range("L134").Value = 11850
$ws.Range("M134").Value = -54120
$ws.Range("N134").Value = -16920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49998.75
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 49998.75
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 49998.75
$ws.Range("N20").Value = -50470.75

$ws.Range("H30").Value = 49998.75
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 49998.75
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 49998.75
$ws.Range("N30").Value = -50180.75

$ws.Range("H128").Value = 49998.75
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49998.75
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49998.75
$ws.Range("N128").Value = -59958.75

$ws.Range("H132").Value = 3400.0557
$ws.Range("I132").Value = 3050.1
$ws.Range("J132").Value = 3837.5
$ws.Range("K132").Value = 9150.299999999999
$ws.Range("L132").Value = 11512.5
$ws.Range("M132").Value = -6620.299999999999
$ws.Range("N132").Value = -16572.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 396.92856
$ws.Range("I5").Value = 433.33334
$ws.Range("J5").Value = 387
$ws.Range("K5").Value = 1300.00002
$ws.Range("L5").Value = 1161
$ws.Range("M5").Value = -1188.00002
$ws.Range("N5").Value = -1385

$ws.Range("H68").Value = 1669.9756
$ws.Range("I68").Value = 763.44446
$ws.Range("J68").Value = 1924.9375
$ws.Range("K68").Value = 2290.33338
$ws.Range("L68").Value = 5774.8125
$ws.Range("M68").Value = -1479.33338
$ws.Range("N68").Value = -7396.8125

$ws.Range("H71").Value = 1669.9756
$ws.Range("I71").Value = 763.44446
$ws.Range("J71").Value = 1924.9375
$ws.Range("K71").Value = 6871.00014
$ws.Range("L71").Value = 17324.4375
$ws.Range("M71").Value = -2815.00014
$ws.Range("N71").Value = -25436.4375

$ws.Range("H92").Value = 854.2
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 854.2
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2562.6
$ws.Range("N92").Value = -5058.6

$ws.Range("H98").Value = 434.875
$ws.Range("I98").Value = 400
$ws.Range("J98").Value = 439.85715
$ws.Range("K98").Value = 1200
$ws.Range("L98").Value = 1319.57145
$ws.Range("M98").Value = 298
$ws.Range("N98").Value = -4315.571449999999

$ws.Range("H131").Value = 6861391
$ws.Range("I131").Value = 166667140
$ws.Range("J131").Value = 12573.028
$ws.Range("K131").Value = 500001420
$ws.Range("L131").Value = 37719.084
$ws.Range("M131").Value = -499996380
$ws.Range("N131").Value = -47799.084

$ws.Range("H132").Value = 1045.421
$ws.Range("I132").Value = 1230
$ws.Range("J132").Value = 1035.1666
$ws.Range("K132").Value = 11070
$ws.Range("L132").Value = 9316.499400000001
$ws.Range("M132").Value = -8540
$ws.Range("N132").Value = -14376.4994

$ws.Range("H133").Value = 31253262
$ws.Range("I133").Value = 62501524
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 187504572
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -187499512
$ws.Range("N133").Value = -25120

$ws.Range("H134").Value = 4802.727
$ws.Range("I134").Value = 4566
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 13698
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -8628
$ws.Range("N134").Value = -25140

$ws.Range("H135").Value = 396.92856
$ws.Range("I135").Value = 433.33334
$ws.Range("J135").Value = 387
$ws.Range("K135").Value = 3900.00006
$ws.Range("L135").Value = 3483
$ws.Range("M135").Value = -1365.00006
$ws.Range("N135").Value = -8553

$ws.Range("H136").Value = 898.6
$ws.Range("I136").Value = 290
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 870
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = 4230
$ws.Range("N136").Value = -20199

$ws.Range("H137").Value = 5436.4546
$ws.Range("I137").Value = 2310
$ws.Range("J137").Value = 5749.1
$ws.Range("K137").Value = 6930
$ws.Range("L137").Value = 17247.3
$ws.Range("M137").Value = -1830
$ws.Range("N137").Value = -27447.3

$ws.Range("H138").Value = 2816
$ws.Range("I138").Value = 2915
$ws.Range("J138").Value = 2750
$ws.Range("K138").Value = 8745
$ws.Range("L138").Value = 8250
$ws.Range("M138").Value = -3605
$ws.Range("N138").Value = -18530

$ws.Range("H139").Value = 36666.668
$ws.Range("I139").Value = 36666.668
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 110000.004
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -104860.004

$ws.Range("H140").Value = 4029.2778
$ws.Range("I140").Value = 793.3333
$ws.Range("J140").Value = 4676.467
$ws.Range("K140").Value = 2379.9999
$ws.Range("L140").Value = 14029.401
$ws.Range("M140").Value = 2800.0001
$ws.Range("N140").Value = -24389.401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14183.546
$ws.Range("I80").Value = 9429.714
$ws.Range("J80").Value = 22502.75
$ws.Range("K80").Value = 9429.714
$ws.Range("L80").Value = 22502.75
$ws.Range("M80").Value = -8431.714
$ws.Range("N80").Value = -24498.75

$ws.Range("H83").Value = 14183.546
$ws.Range("I83").Value = 9429.714
$ws.Range("J83").Value = 22502.75
$ws.Range("K83").Value = 47148.57
$ws.Range("L83").Value = 112513.75
$ws.Range("M83").Value = -42156.57
$ws.Range("N83").Value = -122497.75

$ws.Range("H97").Value = 2056.2727
$ws.Range("I97").Value = 1847
$ws.Range("J97").Value = 2230.6667
$ws.Range("K97").Value = 1847
$ws.Range("L97").Value = 2230.6667
$ws.Range("M97").Value = -1351
$ws.Range("N97").Value = -3222.6667

$ws.Range("H126").Value = 1716802.4
$ws.Range("I126").Value = 2927175.2
$ws.Range("J126").Value = 74153.64
$ws.Range("K126").Value = 8781525.600000001
$ws.Range("L126").Value = 222460.92
$ws.Range("M126").Value = -8779055.600000001
$ws.Range("N126").Value = -227400.92

$ws.Range("H132").Value = 1042080.2
$ws.Range("I132").Value = 1604271.6
$ws.Range("J132").Value = 4188.385
$ws.Range("K132").Value = 4812814.800000001
$ws.Range("L132").Value = 12565.155
$ws.Range("M132").Value = -4810284.800000001
$ws.Range("N132").Value = -17625.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4053.2222
$ws.Range("I7").Value = 2815.1667
$ws.Range("J7").Value = 4672.25
$ws.Range("K7").Value = 2815.1667
$ws.Range("L7").Value = 4672.25
$ws.Range("M7").Value = -2703.1667
$ws.Range("N7").Value = -4896.25

$ws.Range("H100").Value = 2608.5
$ws.Range("I100").Value = 2345
$ws.Range("J100").Value = 4980
$ws.Range("K100").Value = 2345
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -1804
$ws.Range("N100").Value = -6062

$ws.Range("H126").Value = 4053.2222
$ws.Range("I126").Value = 2815.1667
$ws.Range("J126").Value = 4672.25
$ws.Range("K126").Value = 8445.500100000001
$ws.Range("L126").Value = 14016.75
$ws.Range("M126").Value = -5975.500100000001
$ws.Range("N126").Value = -18956.75

$ws.Range("H133").Value = 73163
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 73163
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 73163
$ws.Range("N133").Value = -78223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2446.889
$ws.Range("I81").Value = 2502.75
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 5005.5
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -3944.5
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 2446.889
$ws.Range("I84").Value = 2502.75
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 25027.5
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -19723.5
$ws.Range("N84").Value = -30608

$ws.Range("H96").Value = 3750
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3750
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3750
$ws.Range("N96").Value = -6496

$ws.Range("H100").Value = 241.33333
$ws.Range("I100").Value = 241.33333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 482.66666
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 58.33334000000002

$ws.Range("H126").Value = 4376.44
$ws.Range("I126").Value = 3668.9443
$ws.Range("J126").Value = 6195.7144
$ws.Range("K126").Value = 11006.8329
$ws.Range("L126").Value = 18587.1432
$ws.Range("M126").Value = -8536.832900000001
$ws.Range("N126").Value = -23527.1432

$ws.Range("H132").Value = 2860.3333
$ws.Range("I132").Value = 2418.8
$ws.Range("J132").Value = 3412.25
$ws.Range("K132").Value = 7256.400000000001
$ws.Range("L132").Value = 10236.75
$ws.Range("M132").Value = -4726.400000000001
$ws.Range("N132").Value = -15296.75
